$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.018.12"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.507.26"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.95%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.69"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.51"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +7.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.508.23"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.85%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.78%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.93%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.111.05"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.94%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.28"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.67%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.942.58"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +5.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.503.42"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.35%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.93%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "395.28"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.99%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.30"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.99%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +12.48%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.530"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.09"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.47%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.36"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.48"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +6.52%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.74%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.48"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +7.85%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.57"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.12%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.66"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.904"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +7.03%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.98%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.67"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +7.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.847.72"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.07%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.38"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.48"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.60"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.56%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.53%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +8.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "350.51"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +7.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.09"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.38%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +12.62%  "
